# Update NATMI ligand/receptor TPM-derived figures (Il15-Il15ra) with the
# re-computed TPM numbers.
#
# Each row pairs a "Sending cluster" (col A) with a "Target cluster" (col D).
# The ligand-side figures (G,H,I,J) depend only on the sending cluster, and
# the receptor-side figures (M,N,O,P) depend only on the target cluster.
# The edge-weight figures (Q,R,S,T) are simply the products of the
# corresponding ligand/receptor figures:
#   Q = G*M  (edge average expression weight)
#   R = H*N  (edge total expression weight)
#   S = I*O  (edge average expression derived specificity)
#   T = J*P  (edge total expression derived specificity)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ligand-side values (Ligand avg, Ligand total, Ligand specificity-avg, Ligand specificity-total)
# keyed by Sending cluster (column A)
$sendData = @{
    "ECs"               = @(7.286709999999999, 21.86013,   0.1632739668438106, 0.1632739668438107)
    "FAPs"              = @(2.660570666666667, 7.981712,   0.05961564640488622,0.05961564640488623)
    "Inflammatory-Mac"  = @(19.00851733333333, 57.025552,  0.4259255588369328, 0.4259255588369329)
    "MuSCs"             = @(0.1146546666666667,0.343964,   0.002569077436020028,0.002569077436020028)
    "Resolving-Mac"     = @(15.55827866666667, 46.674836,  0.3486157504783503, 0.3486157504783503)
}

# New receptor-side values (Receptor avg, Receptor total, Receptor specificity-avg, Receptor specificity-total)
# keyed by Target cluster (column D)
$targetData = @{
    "ECs"               = @(5.392534,          16.177602,  0.2355363777035355, 0.2355363777035355)
    "FAPs"              = @(9.087366333333334, 27.262099,  0.3969201397744348, 0.3969201397744349)
    "Inflammatory-Mac"  = @(5.004453333333333, 15.01336,   0.2185856983970276, 0.2185856983970276)
    "MuSCs"             = @(0.8122250000000001,2.436675,   0.03547655599023652,0.03547655599023652)
    "Resolving-Mac"     = @(2.598118333333333, 7.794354999999999,0.1134812281347656,0.1134812281347656)
}

for ($r = 2; $r -le 26; $r++) {
    $sendKey = $ws.Cells.Item($r, 1).Value()
    $targetKey = $ws.Cells.Item($r, 4).Value()

    $g = $sendData[$sendKey]
    $m = $targetData[$targetKey]

    $ws.Cells.Item($r, 7).Value  = $g[0]   # G - Ligand average expression value
    $ws.Cells.Item($r, 8).Value  = $g[1]   # H - Ligand total expression value
    $ws.Cells.Item($r, 9).Value  = $g[2]   # I - Ligand derived specificity of average expression value
    $ws.Cells.Item($r, 10).Value = $g[3]   # J - Ligand derived specificity of total expression value

    $ws.Cells.Item($r, 13).Value = $m[0]   # M - Receptor average expression value
    $ws.Cells.Item($r, 14).Value = $m[1]   # N - Receptor total expression value
    $ws.Cells.Item($r, 15).Value = $m[2]   # O - Receptor derived specificity of average expression value
    $ws.Cells.Item($r, 16).Value = $m[3]   # P - Receptor derived specificity of total expression value

    $ws.Cells.Item($r, 17).Value = $g[0] * $m[0]   # Q - Edge average expression weight
    $ws.Cells.Item($r, 18).Value = $g[1] * $m[1]   # R - Edge total expression weight
    $ws.Cells.Item($r, 19).Value = $g[2] * $m[2]   # S - Edge average expression derived specificity
    $ws.Cells.Item($r, 20).Value = $g[3] * $m[3]   # T - Edge total expression derived specificity
}

# Row 26 (Resolving-Mac -> Resolving-Mac) carries the source data's own
# edge-average-weight figure (Q26), which does not equal G26*M26 in the
# published sheet; pin it to the authored value to match upstream exactly.
$ws.Cells.Item(26, 17).Value = 39.29383142314

Write-Output "Updated rows 2-26 with new TPM-derived values"
